# The presentation currently carries the "Integral" theme (ppt/theme/theme1.xml,
# used by the slide master) and an "Office Theme" (ppt/theme/theme2.xml, used by
# the notes master). The authored change swaps the two theme colour schemes so
# that the deck's main theme becomes the plain "Office Theme" palette (and the
# notes master carries the former "Integral" palette).
#
# The PowerPoint object model only exposes the *active* theme's colour scheme
# (reached the same way from the slide master, the notes master, or any slide),
# so we recolour it here to the target ("Office Theme") palette - the scheme
# that is swapped into ppt/theme/theme1.xml by the authored edit.

$p = $ppt.ActivePresentation
$cs = $p.SlideMaster.Theme.ThemeColorScheme

# Index order matches the OOXML <a:clrScheme> child order:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
# 8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink.
# RGB values use the OLE COLORREF 0x00BBGGRR encoding.

$cs.Colors(1).RGB  = 0        # dk1      #000000
$cs.Colors(2).RGB  = 16777215 # lt1      #FFFFFF
$cs.Colors(3).RGB  = 6968388  # dk2      #44546A
$cs.Colors(4).RGB  = 15132391 # lt2      #E7E6E6
$cs.Colors(5).RGB  = 13998939 # accent1  #5B9BD5
$cs.Colors(6).RGB  = 3243501  # accent2  #ED7D31
$cs.Colors(7).RGB  = 10855845 # accent3  #A5A5A5
$cs.Colors(8).RGB  = 49407    # accent4  #FFC000
$cs.Colors(9).RGB  = 12874308 # accent5  #4472C4
$cs.Colors(10).RGB = 4697456  # accent6  #70AD47
$cs.Colors(11).RGB = 12673797 # hlink    #0563C1
$cs.Colors(12).RGB = 7491477  # folHlink #954F72
